$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.709.93"
$ws.Range("E2").Value = "  +1.28%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.695.90"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.43%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.60"
$ws.Range("E5").Value = "  +1.31%  "

$ws.Range("E6").Value = "  +0.19%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3954"
$ws.Range("E7").Value = "  +0.42%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4065"
$ws.Range("E8").Value = "  +0.48%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.489"
$ws.Range("E9").Value = "  -0.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.001"
$ws.Range("E10").Value = "  +0.48%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.26"
$ws.Range("E11").Value = "  -2.34%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08847"
$ws.Range("E12").Value = "  +0.95%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.256"
$ws.Range("E13").Value = "  -0.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.66"
$ws.Range("E14").Value = "  +1.89%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.069"
$ws.Range("E15").Value = "  +8.46%  "

$ws.Range("E16").Value = "  +0.11%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.702.20"
$ws.Range("E17").Value = "  +0.62%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "100.05"
$ws.Range("E18").Value = "  -0.55%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07031"
$ws.Range("E19").Value = "  -0.62%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.56"
$ws.Range("E20").Value = "  +0.86%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.042"
$ws.Range("E21").Value = "  +4.32%  "

$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.31"
$ws.Range("E23").Value = "  +1.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.708.76"
$ws.Range("E24").Value = "  +1.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.265"
$ws.Range("E25").Value = "  +9.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.362"
$ws.Range("E26").Value = "  +2.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.82"
$ws.Range("E27").Value = "  +1.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.50"
$ws.Range("E28").Value = "  +2.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "136.19"
$ws.Range("E29").Value = "  +1.67%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.188"
$ws.Range("E30").Value = "  +1.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.532"
$ws.Range("E31").Value = "  +1.61%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.885.36"
$ws.Range("E32").Value = "  +0.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.072"
$ws.Range("E33").Value = "  -1.54%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08592"
$ws.Range("E34").Value = "  -1.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.136"
$ws.Range("E35").Value = "  -4.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.42"
$ws.Range("E36").Value = "  +1.00%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2752"
$ws.Range("E37").Value = "  +1.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.908"
$ws.Range("E38").Value = "  -1.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.46"
$ws.Range("E39").Value = "  -1.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09234"
$ws.Range("E40").Value = "  +3.15%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02727"
$ws.Range("E41").Value = "  -1.90%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.473"
$ws.Range("E42").Value = "  +0.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7662"
$ws.Range("E43").Value = "  +1.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.05"
$ws.Range("E44").Value = "  +4.79%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7186"
$ws.Range("E45").Value = "  +0.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.578"
$ws.Range("E46").Value = "  +5.57%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.216"
$ws.Range("E47").Value = "  +1.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.326"
$ws.Range("E49").Value = "  +1.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "139.66"
$ws.Range("E50").Value = "  -0.23%  "

$ws.Range("E51").Value = "  +0.53%  "

